# Auto-generated Excel COM-interop script applying the Typhon_Profits market-data refresh.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) for specific leve rows
# across the ALC, ARM, BSM, CRP, CUL, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 451.6279
$ws.Range("J17").Value = 446.59525
$ws.Range("L17").Value = 1339.78575
$ws.Range("N17").Value = -1675.78575

# Row 19
$ws.Range("H19").Value = 175
$ws.Range("I19").Value = 100
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = 100
$ws.Range("L19").Value = 400
$ws.Range("M19").Value = 75
$ws.Range("N19").Value = -750

# Row 43
$ws.Range("H43").Value = 799.3333
$ws.Range("I43").Value = 798
$ws.Range("K43").Value = 798
$ws.Range("M43").Value = -729

# Row 74
$ws.Range("H74").Value = 2289.4138
$ws.Range("J74").Value = 3665.5557
$ws.Range("L74").Value = 3665.5557
$ws.Range("N74").Value = -5537.5557

# Row 77
$ws.Range("H77").Value = 2289.4138
$ws.Range("J77").Value = 3665.5557
$ws.Range("L77").Value = 18327.7785
$ws.Range("N77").Value = -27687.7785

# Row 112
$ws.Range("H112").Value = 4033315
$ws.Range("J112").Value = 1112.6897
$ws.Range("L112").Value = 3338.0691
$ws.Range("N112").Value = -5554.0691

# Row 132
$ws.Range("H132").Value = 53766.55
$ws.Range("J132").Value = 548
$ws.Range("L132").Value = 1644
$ws.Range("N132").Value = -6704

# Row 138
$ws.Range("H138").Value = 2512.3962
$ws.Range("I138").Value = 484.7647
$ws.Range("J138").Value = 3469.889
$ws.Range("K138").Value = 1454.2941
$ws.Range("L138").Value = 10409.667
$ws.Range("M138").Value = 3685.7059
$ws.Range("N138").Value = -20689.667

$ws = $wb.Worksheets.Item("ARM")
# Row 22
$ws.Range("H22").Value = 1308.6666
$ws.Range("I22").Value = 1308.6666
$ws.Range("K22").Value = 1308.6666
$ws.Range("M22").Value = -1009.6666

# Row 28
$ws.Range("H28").Value = 8860.5
$ws.Range("I28").Value = 5147.3335
$ws.Range("J28").Value = 20000
$ws.Range("K28").Value = 5147.3335
$ws.Range("L28").Value = 20000
$ws.Range("M28").Value = -4955.3335
$ws.Range("N28").Value = -20384

# Row 32
$ws.Range("H32").Value = 18687.838
$ws.Range("I32").Value = 21612.72
$ws.Range("J32").Value = 6500.8335
$ws.Range("K32").Value = 21612.72
$ws.Range("L32").Value = 6500.8335
$ws.Range("M32").Value = -21325.72
$ws.Range("N32").Value = -7074.8335

# Row 41
$ws.Range("H41").Value = 3304
$ws.Range("I41").Value = 3304
$ws.Range("K41").Value = 3304
$ws.Range("M41").Value = -2890

# Row 99
$ws.Range("H99").Value = 8860.5
$ws.Range("I99").Value = 5147.3335
$ws.Range("J99").Value = 20000
$ws.Range("K99").Value = 5147.3335
$ws.Range("L99").Value = 20000
$ws.Range("M99").Value = -2152.3335
$ws.Range("N99").Value = -25990

# Row 105
$ws.Range("H105").Value = 37696
$ws.Range("J105").Value = 37696
$ws.Range("L105").Value = 37696
$ws.Range("N105").Value = -44684

$ws = $wb.Worksheets.Item("BSM")
# Row 8
$ws.Range("H8").Value = 400
$ws.Range("I8").Value = 400
$ws.Range("K8").Value = 400
$ws.Range("M8").Value = -260

# Row 107
$ws.Range("H107").Value = 748.24
$ws.Range("I107").Value = 768.8570999999999
$ws.Range("J107").Value = 640
$ws.Range("K107").Value = 768.8570999999999
$ws.Range("L107").Value = 640
$ws.Range("M107").Value = 1151.1429
$ws.Range("N107").Value = -4480

# Row 134
$ws.Range("H134").Value = 28492.816
$ws.Range("I134").Value = 30061.889
$ws.Range("J134").Value = 249.5
$ws.Range("K134").Value = 90185.667
$ws.Range("L134").Value = 748.5
$ws.Range("M134").Value = -87650.667
$ws.Range("N134").Value = -5818.5

$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 63563356
$ws.Range("I6").Value = 10345016
$ws.Range("J6").Value = 170000030
$ws.Range("K6").Value = 10345016
$ws.Range("L6").Value = 170000030
$ws.Range("M6").Value = -10344903
$ws.Range("N6").Value = -170000256

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 91
$ws.Range("I12").Value = 12.6
$ws.Range("J12").Value = 121.15385
$ws.Range("K12").Value = 37.8
$ws.Range("L12").Value = 363.46155
$ws.Range("M12").Value = 135.2
$ws.Range("N12").Value = -709.46155

# Row 64
$ws.Range("H64").Value = 1798.3334
$ws.Range("I64").Value = 1358
$ws.Range("K64").Value = 4074
$ws.Range("M64").Value = -3804

# Row 67
$ws.Range("H67").Value = 1798.3334
$ws.Range("I67").Value = 1358
$ws.Range("K67").Value = 4074
$ws.Range("M67").Value = -3138

# Row 75
$ws.Range("H75").Value = 3206.077
$ws.Range("I75").Value = 1124.75
$ws.Range("J75").Value = 4131.1113
$ws.Range("K75").Value = 3374.25
$ws.Range("L75").Value = 12393.3339
$ws.Range("M75").Value = -2376.25
$ws.Range("N75").Value = -14389.3339

# Row 78
$ws.Range("H78").Value = 3206.077
$ws.Range("I78").Value = 1124.75
$ws.Range("J78").Value = 4131.1113
$ws.Range("K78").Value = 10122.75
$ws.Range("L78").Value = 37180.00169999999
$ws.Range("M78").Value = -5130.75
$ws.Range("N78").Value = -47164.00169999999

# Row 86
$ws.Range("H86").Value = 846.75
$ws.Range("I86").Value = 858.6667
$ws.Range("J86").Value = 811
$ws.Range("K86").Value = 2576.0001
$ws.Range("L86").Value = 2433
$ws.Range("M86").Value = -1390.0001
$ws.Range("N86").Value = -4805

# Row 87
$ws.Range("H87").Value = 17741.25
$ws.Range("I87").Value = 6857.5
$ws.Range("J87").Value = 28625
$ws.Range("K87").Value = 20572.5
$ws.Range("L87").Value = 85875
$ws.Range("M87").Value = -19324.5
$ws.Range("N87").Value = -88371

# Row 89
$ws.Range("H89").Value = 846.75
$ws.Range("I89").Value = 858.6667
$ws.Range("J89").Value = 811
$ws.Range("K89").Value = 7728.0003
$ws.Range("L89").Value = 7299
$ws.Range("M89").Value = -1800.0003
$ws.Range("N89").Value = -19155

# Row 90
$ws.Range("H90").Value = 17741.25
$ws.Range("I90").Value = 6857.5
$ws.Range("J90").Value = 28625
$ws.Range("K90").Value = 61717.5
$ws.Range("L90").Value = 257625
$ws.Range("M90").Value = -55477.5
$ws.Range("N90").Value = -270105

# Row 93
$ws.Range("H93").Value = 2604.125
$ws.Range("J93").Value = 2500
$ws.Range("L93").Value = 7500
$ws.Range("N93").Value = -11244

# Row 107
$ws.Range("H107").Value = 3969.2666
$ws.Range("J107").Value = 750.1177
$ws.Range("L107").Value = 2250.3531
$ws.Range("N107").Value = -6090.3531

# Row 114
$ws.Range("H114").Value = 2328.5715
$ws.Range("I114").Value = 2325
$ws.Range("J114").Value = 2333.3333
$ws.Range("K114").Value = 6975
$ws.Range("L114").Value = 6999.999899999999
$ws.Range("M114").Value = -3721
$ws.Range("N114").Value = -13507.9999

# Row 129
$ws.Range("H129").Value = 1613.4286
$ws.Range("I129").Value = 571
$ws.Range("J129").Value = 2655.8572
$ws.Range("K129").Value = 1713
$ws.Range("L129").Value = 7967.571599999999
$ws.Range("M129").Value = 3287
$ws.Range("N129").Value = -17967.5716

# Row 131
$ws.Range("H131").Value = 135971.8
$ws.Range("I131").Value = 882.3333
$ws.Range("J131").Value = 147891.45
$ws.Range("K131").Value = 2646.9999
$ws.Range("L131").Value = 443674.35
$ws.Range("M131").Value = 2393.0001
$ws.Range("N131").Value = -453754.35

# Row 137
$ws.Range("H137").Value = 7958.2104
$ws.Range("I137").Value = 25407.25
$ws.Range("J137").Value = 3305.1333
$ws.Range("K137").Value = 76221.75
$ws.Range("L137").Value = 9915.3999
$ws.Range("M137").Value = -71121.75
$ws.Range("N137").Value = -20115.3999

$ws = $wb.Worksheets.Item("LTW")
# Row 20
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

# Row 100
$ws.Range("H100").Value = 6760.6
$ws.Range("I100").Value = 1267.6666
$ws.Range("J100").Value = 15000
$ws.Range("K100").Value = 1267.6666
$ws.Range("L100").Value = 15000
$ws.Range("M100").Value = -726.6666
$ws.Range("N100").Value = -16082

# Row 110
$ws.Range("H110").Value = 2030379.8
$ws.Range("J110").Value = 2030379.8
$ws.Range("L110").Value = 2030379.8
$ws.Range("N110").Value = -2038559.8

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 2457361.8
$ws.Range("I113").Value = 416.77777
$ws.Range("K113").Value = 1250.33331
$ws.Range("M113").Value = 919.66669

# Row 126
$ws.Range("H126").Value = 824.8333
$ws.Range("I126").Value = 824.8333
$ws.Range("K126").Value = 2474.4999
$ws.Range("M126").Value = -4.499899999999798

# Row 132
$ws.Range("H132").Value = 1491.825
$ws.Range("I132").Value = 1147.5807
$ws.Range("J132").Value = 2677.5557
$ws.Range("K132").Value = 3442.7421
$ws.Range("L132").Value = 8032.6671
$ws.Range("M132").Value = -912.7420999999999
$ws.Range("N132").Value = -13092.6671
